# Updated Gmail username and password
# - clears leftover scratch values in B2/C2/B3/C3
# - relabels the existing "Gmail" row as "Gmail 1"
# - adds a new "Gmail 2" row re-using the Gmail username/password
# - moves the active selection to F8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove stray scratch data that was never meant to stay in B2/C2/B3/C3
$ws.Range("B2").ClearContents()
$ws.Range("C2").ClearContents()
$ws.Range("B3").ClearContents()
$ws.Range("C3").ClearContents()

# Existing Gmail account becomes "Gmail 1"
$ws.Range("A3").Value = "Gmail 1"

# New "Gmail 2" row, sharing the same username/password pairing
$ws.Range("A4").Value = "Gmail 2"
$ws.Range("D4").Value = "sandeep.shankar1991@gmail.com"
$ws.Range("E4").Value = "sandeep193"

$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:sandeep.shankar1991@gmail.com") | Out-Null
$ws.Range("D4").Style = "Hyperlink"

$ws.Range("F8").Select() | Out-Null
